# Commit: "added on 19sep 2017"
# Adds a prefix/suffix header pair and a small "today" stamp block to
# Sheet1 (columns F:I, row 5-6), matching the new shared strings /
# number formats / cellXfs that ship with this revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- new headers (row 5) ----------------------------------------------
$ws.Range("F5").Value = "Sh_Prifix"
$ws.Range("G5").Value = "Sh_Sufix"

# --- "today" stamp block (row 6) ---------------------------------------
# Apply the number formats before writing the formula/value so the
# workbook doesn't first auto-detect a generic date format and then
# replace it (that would leave an orphan numFmt/cellXfs entry behind).
$ws.Range("G6").NumberFormat = "mmddyyyy"
$ws.Range("G6").Formula = "=TODAY()"

$ws.Range("H6").NumberFormat = "0_);(0)"
$ws.Range("I6").NumberFormat = "mmddyyyy"

# F6 is a quote-prefixed text "01" formatted with a plain "0" number
# format (mirrors how the source cell was entered: '01).
$ws.Range("F6").Value = "'01"
$ws.Range("F6").NumberFormat = "0"

# --- column sizing for the new block ------------------------------------
$ws.Columns.Item(8).ColumnWidth = 15.28515625
$ws.Columns.Item(9).ColumnWidth = 9.7109375

# --- print setup ---------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- leave the selection where the author left it -----------------------
$ws.Range("F6").Select()
